$wb = $excel.ActiveWorkbook

# Insert a new worksheet named "CalcLog" right after "UserPermission"
# (i.e. before "KaikeiInf"), mirroring Worksheets.Add(Before:=...) semantics.
$afterSheet = $wb.Worksheets.Item("UserPermission")
$ws = $wb.Worksheets.Add($null, $afterSheet)
$ws.Name = "CalcLog"

# ---- Header row (row 1) ----
$ws.Cells.Item(1, 1).Value = "hp_id"
$ws.Cells.Item(1, 2).Value = "pt_id"
$ws.Cells.Item(1, 3).Value = "raiin_no"
$ws.Cells.Item(1, 4).Value = "seq_no"
$ws.Cells.Item(1, 5).Value = "sin_date"
$ws.Cells.Item(1, 6).Value = "log_sbt"
$ws.Cells.Item(1, 7).Value = "text"
$ws.Cells.Item(1, 8).Value = "create_date"
$ws.Cells.Item(1, 9).Value = "create_id"
$ws.Cells.Item(1, 10).Value = "create_machine"
$ws.Cells.Item(1, 11).Value = "update_date"
$ws.Cells.Item(1, 12).Value = "update_id"
$ws.Cells.Item(1, 13).Value = "update_machine"
$ws.Cells.Item(1, 14).Value = "del_item_cd"
$ws.Cells.Item(1, 15).Value = "del_sbt"
$ws.Cells.Item(1, 16).Value = "is_warning"
$ws.Cells.Item(1, 17).Value = "hoken_id"
$ws.Cells.Item(1, 18).Value = "item_cd"
$ws.Cells.Item(1, 19).Value = "term_cnt"
$ws.Cells.Item(1, 20).Value = "term_sbt"

# ---- Data row (row 2) ----
$ws.Cells.Item(2, 1).Value = 998
$ws.Cells.Item(2, 2).Value = 12345
$ws.Cells.Item(2, 3).Value = 1234321
$ws.Cells.Item(2, 4).Value = 1
$ws.Cells.Item(2, 5).Value = 20180807
$ws.Cells.Item(2, 6).Value = 2

# Write the "SmartKarte" cells before the quote-prefixed "test" cell so the
# shared-string table gets the same append order as the target workbook.
$ws.Cells.Item(2, 10).Value = "SmartKarte"
$ws.Cells.Item(2, 13).Value = "SmartKarte"

$ws.Cells.Item(2, 7).Value = "'test"

$ws.Cells.Item(2, 8).Value = 45044.845445567131
$ws.Cells.Item(2, 8).NumberFormat = "mm:ss.0"
$ws.Cells.Item(2, 9).Value = 2

$ws.Cells.Item(2, 11).Value = 45044.845445567131
$ws.Cells.Item(2, 11).NumberFormat = "mm:ss.0"
$ws.Cells.Item(2, 12).Value = 2

$ws.Cells.Item(2, 15).Value = 0
$ws.Cells.Item(2, 16).Value = 0
$ws.Cells.Item(2, 17).Value = 0
$ws.Cells.Item(2, 19).Value = 0
$ws.Cells.Item(2, 20).Value = 0

# ---- Column width (column J / index 10) ----
$ws.Columns.Item(10).ColumnWidth = 14.88671875

# ---- Selection / view state ----
$ws.Range("F7").Select()
